# inputs_outputs sheet: point the "current forecast basic" tool at this
# machine's paths and flip the "new layer" flag on, per
# status_exists_for_control.py (lines 3-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: description + location of the basic-forecast tool itself.
$ws.Range("A2").Value = "מיקום תוכנת תחזית בסיס"
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Color = 0

$ws.Range("B2").Value = "C:\Users\dpere\Documents\JTMT\forecast\create_forecast_basic\current"

# Row 3: output-by-version location.
$ws.Range("B3").Value = "C:\Users\dpere\Documents\JTMT\Projects\תחזיות_דמוגרפיות\קבצי עבודה\142_מתחם_אנגל\בהת\For_approval\Reference_tabels\shp"

# Row 4: "new layer" flag, False -> True. Stored as literal text (matches
# the existing "TRUE"/"FALSE" custom number format), so force text entry.
$ws.Range("B4").Value = "'True"

# Row 5: new-layer location, previously blank.
$ws.Range("B5").Value = "C:\Users\dpere\Documents\JTMT\Projects\תחזיות_דמוגרפיות\קבצי עבודה\142_מתחם_אנגל\בהת\For_approval\Reference_tabels\shp\TAZ_V4_230518_Published.shp"

# Final selection moves from B4 to A2.
$ws.Range("A2").Select() | Out-Null
